$d = $word.ActiveDocument

# 1) "3.11." -> "4.11. " (trailing space preserved)
$d.Content.Find.Execute("3.11.", $true, $false, $false, $false, $false, $true, 1, $false, "4.11. ", 2)

# 2) Fill the (previously empty) second paragraph with the "SQL Befehle" sentence,
#    reproducing the proofing-error markup / run split Word's own spell & grammar
#    checker would have inserted while the text was typed.
$p2 = $d.Paragraphs(2)
$xmlSql = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>SQL Befehle</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> in Word geschrieben. Jetzt f&#252;ge ich sie in Python ein.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xmlSql)

# 3) Insert a brand-new third paragraph listing the Python classes that were created.
$p3 = $d.Paragraphs(3)
$xmlClasses = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Erst einmal die Klassen </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>User,</w:t></w:r><w:r><w:t>RRule</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>VEvent</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> und </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>VCalendar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> erstellt.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($xmlClasses)

# 4) The original trailing empty paragraph got pushed down to position 4 by the
#    InsertXML calls above; remove it (together with the paragraph mark that now
#    separates it from paragraph 3) so the document ends with exactly 3 paragraphs.
$pCount = $d.Paragraphs.Count
$pNew3 = $d.Paragraphs(3)
$pLast = $d.Paragraphs($pCount)
$trailing = $d.Range($pNew3.Range.End - 1, $pLast.Range.End)
$trailing.Delete()
